# Fixed wording and completed the first draft of use cases
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shorten / simplify the "two flight numbers" requirement (row 31 / C31)
$ws.Range("C31").Value = "Það þarf að vera hægt að skrá tvö flugnúmer á vinnuferð "

# Fix "að að" typo -> "að" (row 36 / C36)
$ws.Range("C36").Value = "Það þarf að hægt að breyta nafni á tengiliði og/eða neyðarsímanúmeri sem skráð er fyrir áfangastað"

# Fix wording "flugvélategund" -> "flugvélategundir" (row 41 / C41)
$ws.Range("C41").Value = "Það þarf að vera hægt að sjá lista yfir allar flugvélategundir og hve margir flugmenn hafa réttindi á viðkomandi tegund"

# Update the current selection on the "kröfur" sheet to reflect where the author left off
$ws.Activate()
$ws.Range("L20").Select()
